$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "AUTO_OUTLF76E1"
$ws.Range("Q2").Value = 24.80217461268756
$ws.Range("AY2").Value = "'03973832500"

# Row 3
$ws.Range("B3").Value = "AUTO_OUTL5C0F0"
$ws.Range("Q3").Value = 24.661119311806317
$ws.Range("AY3").Value = "'03976278100"

# Row 4
$ws.Range("B4").Value = "AUTO_OUTL721D0"
$ws.Range("Q4").Value = 24.444845834022193
$ws.Range("AY4").Value = "'03978246500"
